# ExcelUtil readData test-fixture update:
#  - B3:B12 get distinct password strings (previously all shared "dinn2535xfvhjmmbk")
#  - B2 keeps its original text but is re-pointed to a fresh shared-string slot
#  - the 10 separate per-cell hyperlinks that used to live on B3..B12 collapse
#    into a single hyperlink covering the B3:B12 range
#  - the active selection moves from G6 to C8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The engine only supports wiping the *entire* worksheet hyperlink collection
# in one shot (per-item .Delete() is a no-op here), so drop them all and
# rebuild the 13 we need, in the same relative order as the original 22 so
# the r:id numbering comes out sequential and tidy.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:dineshkumar.icon@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Dinnu@247", "", "", "Dinnu@247") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:dineshkumar.icon.dk@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:dineshkumar.icon@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:dineshkumar.icon@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:dineshkumar.icon@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "mailto:dineshkumar.icon@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A12"), "mailto:dineshkumar.icon@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:dineshkumar.icon.dk@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:dineshkumar.icon.dk@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:dineshkumar.icon.dk@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A11"), "mailto:dineshkumar.icon.dk@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3:B12"), "mailto:Dinnu@247", "", "", "Dinnu@247") | Out-Null

# Hyperlinks.Add(..., TextToDisplay:"Dinnu@247") stamps that literal text into
# the cells it touches, so set the real per-row values afterwards.
$ws.Range("B2").Value = "adsdfhghk52346421"
$ws.Range("B3").Value = "adsdfhghk52346422"
$ws.Range("B4").Value = "adsdfhghk52346423"
$ws.Range("B5").Value = "adsdfhghk52346424"
$ws.Range("B6").Value = "adsdfhghk52346425"
$ws.Range("B7").Value = "adsdfhghk52346426"
$ws.Range("B8").Value = "adsdfhghk52346427"
$ws.Range("B9").Value = "adsdfhghk52346428"
$ws.Range("B10").Value = "adsdfhghk52346429"
$ws.Range("B11").Value = "adsdfhghk52346430"
$ws.Range("B12").Value = "adsdfhghk52346431"

# Hyperlinks.Add() also silently re-keys the style of any cell it touches;
# put the original "Hyperlink" cell style back on all of them.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("A5").Style = "Hyperlink"
$ws.Range("A6").Style = "Hyperlink"
$ws.Range("A7").Style = "Hyperlink"
$ws.Range("A8").Style = "Hyperlink"
$ws.Range("A9").Style = "Hyperlink"
$ws.Range("A10").Style = "Hyperlink"
$ws.Range("A11").Style = "Hyperlink"
$ws.Range("A12").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3:B12").Style = "Hyperlink"

# Move the saved selection from G6 to C8.
$ws.Range("C8").Select() | Out-Null
